# "cores nos graficos de barras" - the table behind the bar chart is
# re-sorted alphabetically by student name (column A) instead of by unit
# (column B), which is why every unit's rows (and thus the chart's per-row
# colors) land in a new order. Replicate the re-sort here: the header row
# (row 1) is left untouched, and A2:C32 is reordered by column A ascending.
# xlAscending = 1, xlNo (no header in the selected range) = 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:C32")
$keyRange  = $ws.Range("A2:A32")

$dataRange.Sort($keyRange, 1, $null, $null, 1, $null, $null, 2)
